# Update "想去人数" (want-to-go count) figures that were refreshed by the
# site's scheduled data pull. Applies to sheet "展览" (F2,F3,F11,F13,F16,
# F18,F21,F24,F26,F28,F30,F31,F32,F33,F35,F36,F38), sheet "本地生活" (F4),
# and sheet "全部类型" (F4,F5,F6,F11,F13,F16,F18,F21,F27,F29,F31,F36,F37,
# F38,F40,F47).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 12806
$ws1.Range("F3").Value = 7164
$ws1.Range("F11").Value = 146
$ws1.Range("F13").Value = 1015
$ws1.Range("F16").Value = 1018
$ws1.Range("F18").Value = 248
$ws1.Range("F21").Value = 278
$ws1.Range("F24").Value = 164
$ws1.Range("F26").Value = 5235
$ws1.Range("F28").Value = 1429
$ws1.Range("F30").Value = 1365
$ws1.Range("F31").Value = 63
$ws1.Range("F32").Value = 40
$ws1.Range("F33").Value = 1364
$ws1.Range("F35").Value = 5
$ws1.Range("F36").Value = 595
$ws1.Range("F38").Value = 3733

$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F4").Value = 2017

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 2017
$ws4.Range("F5").Value = 12806
$ws4.Range("F6").Value = 7164
$ws4.Range("F11").Value = 146
$ws4.Range("F13").Value = 1015
$ws4.Range("F16").Value = 1018
$ws4.Range("F18").Value = 248
$ws4.Range("F21").Value = 278
$ws4.Range("F27").Value = 164
$ws4.Range("F29").Value = 5235
$ws4.Range("F31").Value = 1429
$ws4.Range("F36").Value = 1365
$ws4.Range("F37").Value = 63
$ws4.Range("F38").Value = 1364
$ws4.Range("F40").Value = 595
$ws4.Range("F47").Value = 3733
